$d = $word.ActiveDocument

# "I use Fedora (as of May 2019, I use Fedora 28)." -> "... September 2019, I use Fedora 29)."
$d.Content.Find.Execute("May 2019, I use Fedora 28).", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "September 2019, I use Fedora 29).", 2)

# "My Arduino projects use C++ 14 ..." -> "... C++ 17 ..."
$d.Content.Find.Execute("C++ 14", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "C++ 17", 2)

# Merge the two runs of the "Fedora 28 ... One main advantage ..." paragraph into one
# (no textual change; replacing with identical text makes the host coalesce the runs)
$d.Content.Find.Execute( `
    "The Linux setup I describe here is based on Fedora 28 but should be easily adaptable to other distributions. One main advantage of Fedora is its repositories include all AVR build tools.", `
    $true, $false, $false, $false, $false, `
    $true, 1, $false, `
    "The Linux setup I describe here is based on Fedora 28 but should be easily adaptable to other distributions. One main advantage of Fedora is its repositories include all AVR build tools.", `
    2)

# avr-g++ / avr-gcc version bump: "7.4" -> "9.2" (replaces both table entries)
$d.Content.Find.Execute("7.4", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "9.2", 2)

# Merge the three runs of "To install it on your system, just type the command line:" into one
$d.Content.Find.Execute( `
    "To install it on your system, just type the command line:", `
    $true, $false, $false, $false, $false, `
    $true, 1, $false, `
    "To install it on your system, just type the command line:", `
    2)
